$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @(3, 1, 0.04027, 0.12081, 0.01318991723029425, 0.01318991723029425, 2, 0.6666666666666666, 0.09554499999999999, 0.286635, 0.01029975823317688, 0.01029975823317688, 0.00384759715, 0.03462837435, 0.0001358529585876448, 0.0001358529585876448)
  3 = @(3, 1, 0.04027, 0.12081, 0.01318991723029425, 0.01318991723029425, 3, 1, 3.721182333333333, 11.163547, 0.4011437372432085, 0.4011437372432086, 0.1498520125633333, 1.34866811307, 0.005291052691688826, 0.005291052691688826)
  4 = @(3, 1, 0.04027, 0.12081, 0.01318991723029425, 0.01318991723029425, 3, 1, 5.459703999999999, 16.379112, 0.5885565045236145, 0.5885565045236146, 0.21986228008, 1.97876052072, 0.007763011580017781, 0.007763011580017781)
  5 = @(3, 1, 0.9943730000000001, 2.983119, 0.3256940046198011, 0.325694004619801, 2, 0.6666666666666666, 0.09554499999999999, 0.286635, 0.01029975823317688, 0.01029975823317688, 0.09500736828499999, 0.855066314565, 0.003354569505579145, 0.003354569505579145)
  6 = @(3, 1, 0.9943730000000001, 2.983119, 0.3256940046198011, 0.325694004619801, 3, 1, 3.721182333333333, 11.163547, 0.4011437372432085, 0.4011437372432086, 3.700243240343667, 33.302189163093, 0.1306501102108938, 0.1306501102108938)
  7 = @(3, 1, 0.9943730000000001, 2.983119, 0.3256940046198011, 0.325694004619801, 3, 1, 5.459703999999999, 16.379112, 0.5885565045236145, 0.5885565045236146, 5.428982245592, 48.860840210328, 0.1916893249033281, 0.1916893249033281)
  8 = @(3, 1, 2.018446666666666, 6.055339999999999, 0.6611160781499047, 0.6611160781499047, 2, 0.6666666666666666, 0.09554499999999999, 0.286635, 0.01029975823317688, 0.01029975823317688, 0.1928524867666666, 1.7356723809, 0.006809335769010091, 0.006809335769010091)
  9 = @(3, 1, 2.018446666666666, 6.055339999999999, 0.6611160781499047, 0.6611160781499047, 3, 1, 3.721182333333333, 11.163547, 0.4011437372432085, 0.4011437372432086, 7.511008076775553, 67.59907269098, 0.2652025743406259, 0.2652025743406259)
  10 = @(3, 1, 2.018446666666666, 6.055339999999999, 0.6611160781499047, 0.6611160781499047, 3, 1, 5.459703999999999, 16.379112, 0.5885565045236145, 0.5885565045236146, 11.02012133978666, 99.18109205807998, 0.3891041680402686, 0.3891041680402688)
}

$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in $data.Keys) {
  $vals = $data[$r]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $addr = "$($cols[$i])$r"
    $ws.Range($addr).Value = $vals[$i]
  }
}
